# Daily attendance processing - 2025-11-28 15:27:02
# Reorders the "Recorded By" list in column G so that the literal token
# "System" (exact case) is moved to the front of the comma-separated
# list, preserving the relative order of the remaining tokens
# (including a lower-case "system" entry, which is left in place
# relative to the other non-"System" tokens).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ","
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $systemParts = @()
    $otherParts = @()
    foreach ($p in $trimmed) {
        if ($p.CompareTo("System") -eq 0) {
            $systemParts += $p
        } else {
            $otherParts += $p
        }
    }

    if ($systemParts.Count -eq 0) { continue }

    $newParts = $systemParts + $otherParts
    $newVal = [string]::Join(", ", $newParts)

    if ($newVal -ne $val) {
        $cell.Value = $newVal
    }
}
